$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values for rows 2-79 were updated from 45190 to 45192
$ws.Range("C2:C79").Value = 45192
